$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header cells from Polish to English
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("E1").Value = "Appearances"
$ws.Range("F1").Value = "1 squad"

# Adjust column widths to match the final layout
# (values chosen so the resulting stored width is as close as possible
# to the target widths of 22.42578125 / 12.5703125 / 19)
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(5).ColumnWidth = 11.666666666666666
$ws.Columns.Item(6).ColumnWidth = 18.166666666666668

# Update the selected cell to match the saved state
$ws.Range("F12").Select()
